# Apply the StructureDefinition metadata refresh:
#  - Version 5.0.0 -> 6.0.0
#  - Date refreshed
#  - Publisher gains a value ("Alvearie Team")
#  - "Contact" / "No display for ContactDetail" row replaced by
#    "Jurisdiction" / "United States of America"
#  - the old duplicate "Contact" row is removed entirely (table shrinks
#    from 21 to 20 rows)
#  - the Elements sheet's root "Extension" element row gets its
#    Short/Definition overwritten with the managed-care-plan-type text

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------

# Version
$wsMeta.Range("B3").Value = "6.0.0"

# Date
$wsMeta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$wsMeta.Range("B9").Value = "Alvearie Team"

# Row 10 used to be "Contact" / "No display for ContactDetail"; it becomes
# "Jurisdiction" / "United States of America"
$wsMeta.Range("A10").Value = "Jurisdiction"
$wsMeta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row;
# delete it outright so everything below shifts up one row
$wsMeta.Rows.Item(11).Delete()

# --- Elements sheet ---------------------------------------------------

# Row 2 is the root "Extension" element; its Short (K) and Definition (L)
# get overwritten with the plan-type name/description
$wsElem.Range("K2").Value = "Managed Care Plan Type"
$wsElem.Range("L2").Value = "Code for the managed care plan type under which the eligible individual is enrolled"
